$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column G with its header, copying the header style used by the
# other header cells (F1) and the new "备注" shared string value.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "备注"

# Size the new column to match the other template columns.
$ws.Columns.Item(7).ColumnWidth = 29.857142857142858

# Add a dropdown list data validation on column A restricting entries to
# "零部件" (component) or "原材料" (raw material).
$ws.Range("A1:A1048576").Validation.Add(3, 1, 1, '"零部件,原材料"')

# Restore the active selection to match the authored workbook.
[void]$ws.Range("F10").Select()
